# Update "想去人数" (interest count) values in column F across the
# relevant worksheets, per the upstream data refresh (gh-pages output
# regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4006
$ws1.Range("F4").Value = 2352
$ws1.Range("F5").Value = 467
$ws1.Range("F8").Value = 20
$ws1.Range("F11").Value = 64
$ws1.Range("F13").Value = 1493
$ws1.Range("F15").Value = 2797

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 39

# --- Sheet "全部类型" (All types, combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4006
$ws4.Range("F4").Value = 2352
$ws4.Range("F5").Value = 467
$ws4.Range("F8").Value = 20
$ws4.Range("F9").Value = 39
$ws4.Range("F12").Value = 64
$ws4.Range("F16").Value = 1493
$ws4.Range("F18").Value = 2797
